# repull data, push all data, mean calculation
# Update column F ("dSF") values for the luzardo_jesús appearance log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = -2
    5  = -1
    6  = -2
    7  = 1
    8  = 2
    9  = -5
    10 = -2
    11 = 9
    12 = -4
    13 = 2
    14 = -1
    15 = 2
    17 = -5
    18 = 2
    19 = 4
    21 = -6
    22 = -2
    23 = 1
    24 = -5
    26 = -3
    28 = -5
    29 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
